$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 40497
$ws.Cells.Item(3, 1).Value = 40678
$ws.Cells.Item(4, 1).Value = 40862
$ws.Cells.Item(5, 1).Value = 41044
$ws.Cells.Item(6, 1).Value = 41228
$ws.Cells.Item(7, 1).Value = 41409
$ws.Cells.Item(8, 1).Value = 41593
$ws.Cells.Item(9, 1).Value = 41774
$ws.Cells.Item(10, 1).Value = 41958
$ws.Cells.Item(11, 1).Value = 42139
$ws.Cells.Item(12, 1).Value = 42323
$ws.Cells.Item(13, 1).Value = 42505
$ws.Cells.Item(14, 1).Value = 42689
$ws.Cells.Item(15, 1).Value = 42870
$ws.Cells.Item(16, 1).Value = 43054
$ws.Cells.Item(17, 1).Value = 43146
$ws.Cells.Item(18, 1).Value = 43235
$ws.Cells.Item(19, 1).Value = 43327
$ws.Cells.Item(20, 1).Value = 43419
$ws.Cells.Item(21, 1).Value = 43511
$ws.Cells.Item(22, 1).Value = 43600
$ws.Cells.Item(23, 1).Value = 43692
$ws.Cells.Item(24, 1).Value = 43784
$ws.Cells.Item(25, 1).Value = 43876
$ws.Cells.Item(26, 1).Value = 43966
$ws.Cells.Item(27, 1).Value = 44058
$ws.Cells.Item(28, 1).Value = 44150
$ws.Cells.Item(29, 1).Value = 44242
$ws.Cells.Item(30, 1).Value = 44331
$ws.Cells.Item(31, 1).Value = 44423
$ws.Cells.Item(32, 1).Value = 44515
$ws.Cells.Item(33, 1).Value = 44607
$ws.Cells.Item(34, 1).Value = 44696
$ws.Cells.Item(35, 1).Value = 44788
$ws.Cells.Item(36, 1).Value = 44880
$ws.Cells.Item(37, 1).Value = 44972
$ws.Cells.Item(38, 1).Value = 45061
$ws.Cells.Item(39, 1).Value = 45153
$ws.Cells.Item(40, 1).Value = 45245
$ws.Cells.Item(41, 1).Value = 45337
$ws.Cells.Item(42, 1).Value = 45427
$ws.Cells.Item(43, 1).Value = 45519
$ws.Cells.Item(44, 1).Value = 45611
$ws.Cells.Item(45, 1).Value = 45703
$ws.Cells.Item(46, 1).Value = 45792
$ws.Cells.Item(47, 1).Value = 45884
